$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 14.321881
$ws.Range("N2").Value = 42.965643
$ws.Range("O2").Value = 0.2949569176783066
$ws.Range("P2").Value = 0.2949569176783066
$ws.Range("Q2").Value = 25.88822732163966
$ws.Range("R2").Value = 232.994045894757
$ws.Range("S2").Value = 0.2949569176783066
$ws.Range("T2").Value = 0.2949569176783066

# Row 3
$ws.Range("N3").Value = 81.25250700000001
$ws.Range("O3").Value = 0.557794259435499
$ws.Range("P3").Value = 0.557794259435499
$ws.Range("Q3").Value = 48.95733485634367
$ws.Range("R3").Value = 440.6160137070931
$ws.Range("S3").Value = 0.557794259435499
$ws.Range("T3").Value = 0.557794259435499

# Row 4
$ws.Range("M4").Value = 7.149790333333333
$ws.Range("N4").Value = 21.449371
$ws.Range("O4").Value = 0.1472488228861944
$ws.Range("P4").Value = 0.1472488228861943
$ws.Range("Q4").Value = 12.92395862326989
$ws.Range("R4").Value = 116.315627609429
$ws.Range("S4").Value = 0.1472488228861944
$ws.Range("T4").Value = 0.1472488228861943
